$d = $word.ActiveDocument

$paragraphsXml = @(
    '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Tworzymy certyfikat</w:t></w:r></w:p>',
    '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Tworzymy bazę danych na serwerze,</w:t></w:r></w:p>',
    '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Wchodzimy w naszą domenę przez </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>sftp</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>bitvise</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:br/><w:t xml:space="preserve">tworzymy folder </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>api</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>',
    '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Za pomocą konsoli </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>bsh</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> w </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>bitvise</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>klonujmey</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> tam </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>backend</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> z gita</w:t></w:r></w:p>',
    '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Z </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>directAdmina</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> tworzymy aplikację </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>nodeJS</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>',
    '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Ustawiamy wszystko, dodajemy prefiks </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>api</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>żebu</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> oddzielnie szły tematy na be i fe</w:t></w:r></w:p>',
    '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Wybieramy plik startowy /</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>dist</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>/index.js</w:t></w:r></w:p>',
    '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Tworzymy aplikację</w:t></w:r></w:p>',
    '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Instalacja </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>npm</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> na serwerze</w:t></w:r></w:p>',
    '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Instalacja </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>js</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> scripts – budowa aplikacji</w:t></w:r></w:p>',
    '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:br/></w:r></w:p>'
)

foreach ($pxml in $paragraphsXml) {
    $lastPara = $d.Paragraphs.Last
    $r = $lastPara.Range
    $r.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Last
    $r2 = $newPara.Range
    $r2.InsertXML($pxml)
}

Write-Output ("Paragraphs inserted; total paragraph count = " + $d.Paragraphs.Count)
